$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..5) {
    $ws.Cells.Item($row, 2).Value = "No"
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 10).Value = 0
}
